# Update loading_percent.xlsx values for "case with 380 kV done"
# (Case_2_205 res_line loading percentages, rows 2-25, columns B,C,E,F,G,I,J,K,M,N)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.95032329936273
$ws.Range("C2").Value = 8.377319678785501
$ws.Range("E2").Value = 15.61354733220693
$ws.Range("F2").Value = 45.5150150595041
$ws.Range("G2").Value = 3.703368190076126
$ws.Range("I2").Value = 27.90804241427117
$ws.Range("J2").Value = 9.355610019263723
$ws.Range("K2").Value = 12.83475567117315
$ws.Range("M2").Value = 18.20853990229797
$ws.Range("N2").Value = 21.85661907271963
$ws.Range("B3").Value = 11.75483818139408
$ws.Range("C3").Value = 8.239845723147175
$ws.Range("E3").Value = 15.55097186909677
$ws.Range("F3").Value = 45.40045695178148
$ws.Range("G3").Value = 3.706015662677202
$ws.Range("I3").Value = 27.9565210270328
$ws.Range("J3").Value = 9.375836252432247
$ws.Range("K3").Value = 12.71316611464694
$ws.Range("M3").Value = 18.13629358692559
$ws.Range("N3").Value = 21.91653272859114
$ws.Range("B4").Value = 11.63647734020331
$ws.Range("C4").Value = 8.156875839012327
$ws.Range("E4").Value = 15.51562652835702
$ws.Range("F4").Value = 45.34072724132225
$ws.Range("G4").Value = 3.707725912344236
$ws.Range("I4").Value = 27.99132133975244
$ws.Range("J4").Value = 9.388919298169538
$ws.Range("K4").Value = 12.64122137246965
$ws.Range("M4").Value = 18.09566728336914
$ws.Range("N4").Value = 21.9552578708792
$ws.Range("B5").Value = 11.58873101128005
$ws.Range("C5").Value = 8.123473766226969
$ws.Range("E5").Value = 15.5020072581182
$ws.Range("F5").Value = 45.31906877630643
$ws.Range("G5").Value = 3.708444222263829
$ws.Range("I5").Value = 28.00676589209519
$ws.Range("J5").Value = 9.394418168760676
$ws.Range("K5").Value = 12.61261608488916
$ws.Range("M5").Value = 18.08006290777535
$ws.Range("N5").Value = 21.97152674918894
$ws.Range("B6").Value = 11.5808341211265
$ws.Range("C6").Value = 8.117953425623435
$ws.Range("E6").Value = 15.4997934489618
$ws.Range("F6").Value = 45.3156347395927
$ws.Range("G6").Value = 3.708564789879885
$ws.Range("I6").Value = 28.00940665414081
$ws.Range("J6").Value = 9.395341376971809
$ws.Range("K6").Value = 12.60791012294423
$ws.Range("M6").Value = 18.07752960191111
$ws.Range("N6").Value = 21.97425768486282
$ws.Range("B7").Value = 11.63583135742688
$ws.Range("C7").Value = 8.156423650256979
$ws.Range("E7").Value = 15.51543966531863
$ws.Range("F7").Value = 45.34042427247335
$ws.Range("G7").Value = 3.707735513103164
$ws.Range("I7").Value = 27.99152451984365
$ws.Range("J7").Value = 9.388992779389902
$ws.Range("K7").Value = 12.64083266583064
$ws.Range("M7").Value = 18.09545297050581
$ws.Range("N7").Value = 21.95547530150753
$ws.Range("B8").Value = 11.88261433968447
$ws.Range("C8").Value = 8.329648088466511
$ws.Range("E8").Value = 15.591339535399
$ws.Range("F8").Value = 45.47332268408775
$ws.Range("G8").Value = 3.704263505146824
$ws.Range("I8").Value = 27.92371137951926
$ws.Range("J8").Value = 9.36244646522635
$ws.Range("K8").Value = 12.79228773315622
$ws.Range("M8").Value = 18.18286233780692
$ws.Range("N8").Value = 21.87687554728262
$ws.Range("B9").Value = 12.37682781303091
$ws.Range("C9").Value = 8.678690887238695
$ws.Range("E9").Value = 15.76407601158204
$ws.Range("F9").Value = 45.81737943993907
$ws.Range("G9").Value = 3.698123553256232
$ws.Range("I9").Value = 27.830791543474
$ws.Range("J9").Value = 9.315638034422546
$ws.Range("K9").Value = 13.10938520671625
$ws.Range("M9").Value = 18.38329647176801
$ws.Range("N9").Value = 21.73808265784492
$ws.Range("B10").Value = 12.74226424641001
$ws.Range("C10").Value = 8.938068331625711
$ws.Range("E10").Value = 15.9048110595362
$ws.Range("F10").Value = 46.11986708759709
$ws.Range("G10").Value = 3.694015478184833
$ws.Range("I10").Value = 27.78709631905667
$ws.Range("J10").Value = 9.284419349777915
$ws.Range("K10").Value = 13.35262774056384
$ws.Range("M10").Value = 18.54737962946579
$ws.Range("N10").Value = 21.64541324315353
$ws.Range("B11").Value = 12.90821221015296
$ws.Range("C11").Value = 9.056131489032362
$ws.Range("E11").Value = 15.97166189538221
$ws.Range("F11").Value = 46.26796786234541
$ws.Range("G11").Value = 3.69223310804351
$ws.Range("I11").Value = 27.7725815994211
$ws.Range("J11").Value = 9.270899810615589
$ws.Range("K11").Value = 13.46507923151145
$ws.Range("M11").Value = 18.6254722891649
$ws.Range("N11").Value = 21.60526480786388
$ws.Range("B12").Value = 12.97094604308329
$ws.Range("C12").Value = 9.100803174513334
$ws.Range("E12").Value = 15.99736754133095
$ws.Range("F12").Value = 46.32552962358621
$ws.Range("G12").Value = 3.691570521162491
$ws.Range("I12").Value = 27.76785815045931
$ws.Range("J12").Value = 9.265877936933144
$ws.Range("K12").Value = 13.50788332528183
$ws.Range("M12").Value = 18.65552130030094
$ws.Range("N12").Value = 21.59034958635084
$ws.Range("B13").Value = 12.95744096244756
$ws.Range("C13").Value = 9.091184662921368
$ws.Range("E13").Value = 15.9918142600118
$ws.Range("F13").Value = 46.31306738868355
$ws.Range("G13").Value = 3.691712672558131
$ws.Range("I13").Value = 27.76884102558859
$ws.Range("J13").Value = 9.266955148958539
$ws.Range("K13").Value = 13.49865546027484
$ws.Range("M13").Value = 18.64902880606184
$ws.Range("N13").Value = 21.59354904022859
$ws.Range("B14").Value = 12.91337587048983
$ws.Range("C14").Value = 9.059807642622546
$ws.Range("E14").Value = 15.97376897969933
$ws.Range("F14").Value = 46.27267403982056
$ws.Range("G14").Value = 3.692178349338281
$ws.Range("I14").Value = 27.77217749755753
$ws.Range("J14").Value = 9.270484702550306
$ws.Range("K14").Value = 13.46859653307862
$ws.Range("M14").Value = 18.62793499181041
$ws.Range("N14").Value = 21.60403195224887
$ws.Range("B15").Value = 12.88636889462462
$ws.Range("C15").Value = 9.04058223740727
$ws.Range("E15").Value = 15.9627661104022
$ws.Range("F15").Value = 46.24812364272702
$ws.Range("G15").Value = 3.692465197173719
$ws.Range("I15").Value = 27.77432189534709
$ws.Range("J15").Value = 9.272659365736651
$ws.Range("K15").Value = 13.45021225749357
$ws.Range("M15").Value = 18.61507596103286
$ws.Range("N15").Value = 21.61049054475292
$ws.Range("B16").Value = 12.73140766997855
$ws.Range("C16").Value = 8.930350044821514
$ws.Range("E16").Value = 15.90049776984918
$ws.Range("F16").Value = 46.11039705553648
$ws.Range("G16").Value = 3.694133693291727
$ws.Range("I16").Value = 27.78815296042146
$ws.Range("J16").Value = 9.285316573706966
$ws.Range("K16").Value = 13.34531176594398
$ws.Range("M16").Value = 18.54234390326434
$ws.Range("N16").Value = 21.64807738455461
$ws.Range("B17").Value = 12.63622227841967
$ws.Range("C17").Value = 8.862710773158645
$ws.Range("E17").Value = 15.86301159354795
$ws.Range("F17").Value = 46.02857414790643
$ws.Range("G17").Value = 3.695179345610723
$ws.Range("I17").Value = 27.79801272385438
$ws.Range("J17").Value = 9.293255759077073
$ws.Range("K17").Value = 13.28139283496724
$ws.Range("M17").Value = 18.49859551184228
$ws.Range("N17").Value = 21.67164943547042
$ws.Range("B18").Value = 12.58145053240279
$ws.Range("C18").Value = 8.823815826779246
$ws.Range("E18").Value = 15.84171790223668
$ws.Range("F18").Value = 45.98250183652589
$ws.Range("G18").Value = 3.6957889148919
$ws.Range("I18").Value = 27.80418844522293
$ws.Range("J18").Value = 9.297886384177756
$ws.Range("K18").Value = 13.24479996995954
$ws.Range("M18").Value = 18.47375856350527
$ws.Range("N18").Value = 21.68539642699866
$ws.Range("B19").Value = 12.56290373066565
$ws.Range("C19").Value = 8.810649719506177
$ws.Range("E19").Value = 15.83455462817224
$ws.Range("F19").Value = 45.96707347162571
$ws.Range("G19").Value = 3.695996704324436
$ws.Range("I19").Value = 27.80636605124779
$ws.Range("J19").Value = 9.299465275280811
$ws.Range("K19").Value = 13.23244087308102
$ws.Range("M19").Value = 18.46540575511879
$ws.Range("N19").Value = 21.69008339482988
$ws.Range("B20").Value = 12.6463578515589
$ws.Range("C20").Value = 8.869910460982128
$ws.Range("E20").Value = 15.86697450512845
$ws.Range("F20").Value = 46.03718208041517
$ws.Range("G20").Value = 3.695067192407336
$ws.Range("I20").Value = 27.79691089191999
$ws.Range("J20").Value = 9.292403975980635
$ws.Range("K20").Value = 13.28817962493555
$ws.Range("M20").Value = 18.50321900231076
$ws.Range("N20").Value = 21.66912059649274
$ws.Range("B21").Value = 12.92632226884928
$ws.Range("C21").Value = 9.069025180490051
$ws.Range("E21").Value = 15.97905884542702
$ws.Range("F21").Value = 46.28449865386893
$ws.Range("G21").Value = 3.692041233928448
$ws.Range("I21").Value = 27.77117650358146
$ws.Range("J21").Value = 9.269445339098557
$ws.Range("K21").Value = 13.47741986608964
$ws.Range("M21").Value = 18.63411796687153
$ws.Range("N21").Value = 21.60094505153832
$ws.Range("B22").Value = 13.10864278784683
$ws.Range("C22").Value = 9.19892643619589
$ws.Range("E22").Value = 16.05458168334922
$ws.Range("F22").Value = 46.45474056352294
$ws.Range("G22").Value = 3.69013559573107
$ws.Range("I22").Value = 27.75886365983033
$ws.Range("J22").Value = 9.255009733796106
$ws.Range("K22").Value = 13.60237096943992
$ws.Range("M22").Value = 18.72243899292238
$ws.Range("N22").Value = 21.55806754174564
$ws.Range("B23").Value = 13.01141546535783
$ws.Range("C23").Value = 9.129631802214739
$ws.Range("E23").Value = 16.01407161351246
$ws.Range("F23").Value = 46.36310242105836
$ws.Range("G23").Value = 3.691146105061897
$ws.Range("I23").Value = 27.76502241204091
$ws.Range("J23").Value = 9.262662333943609
$ws.Range("K23").Value = 13.53557826329507
$ws.Range("M23").Value = 18.67505338192602
$ws.Range("N23").Value = 21.58079860398511
$ws.Range("B24").Value = 12.64177570465843
$ws.Range("C24").Value = 8.866655504169524
$ws.Range("E24").Value = 15.86518206915725
$ws.Range("F24").Value = 46.03328741407796
$ws.Range("G24").Value = 3.695117870668321
$ws.Range("I24").Value = 27.79740745023745
$ws.Range("J24").Value = 9.292788860615273
$ws.Range("K24").Value = 13.28511083154977
$ws.Range("M24").Value = 18.50112774090161
$ws.Range("N24").Value = 21.67026327673884
$ws.Range("B25").Value = 12.24243972666973
$ws.Range("C25").Value = 8.583552189840235
$ws.Range("E25").Value = 15.71486116960264
$ws.Range("F25").Value = 45.71547878671507
$ws.Range("G25").Value = 3.699713474897536
$ws.Range("I25").Value = 27.85162348545098
$ws.Range("J25").Value = 9.327741981594579
$ws.Range("K25").Value = 13.02165169445415
$ws.Range("M25").Value = 18.32605212084123
$ws.Range("N25").Value = 21.77399276076031
